# Edit the hypotheses_predictions workbook per the July commit.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("strategies")
$ws2 = $wb.Worksheets.Item("meta")

# --- Data edits on "strategies" sheet ---
# New "Notes:" section appended below existing notes (rows 28-29)
$ws1.Range("A29").Value = "Dan says: Tolerant-avoider is about other functional traits; opprotunistic-conservative is about phenological traits"

$ws1.Range("A28").Value = "Notes:"
$ws1.Range("A28").Font.Bold = $true

# Row 17 (height row): F/G columns updated
$ws1.Range("F17").Value = "Height: unimportant"
$ws1.Range("G17").Value = "Height unimportant"

# Row 16 (leaf N row): F/G columns updated with new hypothesis text
$ws1.Range("F16").Value = "Lower leaf N (more structure)"
$ws1.Range("G16").Value = "No leaf N prediction"

# --- Data edits on "meta" sheet ---
$ws2.Range("A3").Value = "Updated by Dan some days later"
$ws2.Range("A4").Value = "Then a couple edits by Lizzie and Dan together in early July"

# --- View/selection state ---
$ws1.Range("C25").Select()
$ws2.Range("A5").Select()
$ws2.Activate()
